$wb = $excel.ActiveWorkbook

# Add a new worksheet that will become the first tab ("Sheet1")
$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, [System.Type]::Missing, [System.Type]::Missing, [System.Type]::Missing)
$newSheet.Move($firstSheet)
$newSheet.Name = "Sheet1"

# Fill in the caption / header content
$newSheet.Range("A1").Value = "New Phytologist Supporting Information"
$newSheet.Range("A2").Value = "Photographs as an essential biodiversity resource: drivers of gaps in the vascular plant photographic record"
$newSheet.Range("A3").Value = "Thomas Mesaglio, Hervé Sauquet, David Coleman, Elizabeth Wenk, William K Cornwell"
$newSheet.Range("A4").Value = "Accepted 8 February 2023"
$newSheet.Range("A6").Value = "Caption"
$newSheet.Range("A7").Value = "All changes made to species included in the original Australian Plant Census data download. These changes cover both omissions (i.e., treatment of putative native species as non-native) and corrections of erroneous location data. Explanations are provided for all changes."

# Bold the two heading cells (A1 and A6)
$newSheet.Range("A1").Font.Bold = $true
$newSheet.Range("A6").Font.Bold = $true

# Make the new sheet the active / selected tab
$newSheet.Activate() | Out-Null
$newSheet.Select() | Out-Null
$newSheet.Range("F19").Select() | Out-Null
